$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 7000
$ws.Range("K64").Value = 7000
$ws.Range("M64").Value = -6752

$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 7000
$ws.Range("K67").Value = 7000
$ws.Range("M67").Value = -6142

$ws.Range("H74").Value = 7000
$ws.Range("I74").Value = 7000
$ws.Range("K74").Value = 7000
$ws.Range("M74").Value = -6064

$ws.Range("H77").Value = 7000
$ws.Range("I77").Value = 7000
$ws.Range("K77").Value = 35000
$ws.Range("M77").Value = -30320

$ws.Range("H99").Value = 848.3333
$ws.Range("I99").Value = 704.3333
$ws.Range("J99").Value = 992.3333
$ws.Range("K99").Value = 2112.9999
$ws.Range("L99").Value = 2976.9999
$ws.Range("M99").Value = -614.9998999999998
$ws.Range("N99").Value = -5972.9999

$ws.Range("H111").Value = 1828.1666
$ws.Range("I111").Value = 1828.1666
$ws.Range("K111").Value = 5484.4998
$ws.Range("M111").Value = -2417.4998

$ws.Range("H113").Value = 1948
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 2209.862
$ws.Range("I132").Value = 1276.4231
$ws.Range("K132").Value = 3829.2693
$ws.Range("M132").Value = -1299.2693

$ws.Range("H137").Value = 1559
$ws.Range("I137").Value = 1198.75
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 3596.25
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -1046.25
$ws.Range("N137").Value = -14100

$ws.Range("H138").Value = 2144.0894
$ws.Range("I138").Value = 1915.1428
$ws.Range("K138").Value = 5745.428400000001
$ws.Range("M138").Value = -605.4284000000007

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5527.645
$ws.Range("I32").Value = 3773.88
$ws.Range("K32").Value = 3773.88
$ws.Range("M32").Value = -3486.88

$ws.Range("H74").Value = 11666.333
$ws.Range("I74").Value = 4999
$ws.Range("K74").Value = 4999
$ws.Range("M74").Value = -4125

$ws.Range("H77").Value = 11666.333
$ws.Range("I77").Value = 4999
$ws.Range("K77").Value = 24995
$ws.Range("M77").Value = -20627

$ws.Range("H104").Value = 54500
$ws.Range("J104").Value = 54500
$ws.Range("L104").Value = 54500
$ws.Range("N104").Value = -61488

$ws.Range("H106").Value = 60666.668
$ws.Range("J106").Value = 60666.668
$ws.Range("L106").Value = 60666.668
$ws.Range("N106").Value = -63190.668

$ws.Range("H132").Value = 4499.643
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 10500
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3649
$ws.Range("I20").Value = 2815.4546
$ws.Range("J20").Value = 4958.857
$ws.Range("K20").Value = 2815.4546
$ws.Range("L20").Value = 4958.857
$ws.Range("M20").Value = -2568.4546
$ws.Range("N20").Value = -5452.857

$ws.Range("H99").Value = 2589.8
$ws.Range("I99").Value = 2589.8
$ws.Range("K99").Value = 2589.8
$ws.Range("M99").Value = -1091.8

$ws.Range("H134").Value = 1724.75
$ws.Range("I134").Value = 1724.75
$ws.Range("K134").Value = 5174.25
$ws.Range("M134").Value = -2639.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 11313.538
$ws.Range("I3").Value = 9831
$ws.Range("K3").Value = 9831
$ws.Range("M3").Value = -9718

$ws.Range("H31").Value = 5265.2
$ws.Range("I31").Value = 4564.154
$ws.Range("J31").Value = 6567.143
$ws.Range("K31").Value = 4564.154
$ws.Range("L31").Value = 6567.143
$ws.Range("M31").Value = -4269.154
$ws.Range("N31").Value = -7157.143

$ws.Range("H33").Value = 549.5
$ws.Range("I33").Value = 549.5
$ws.Range("K33").Value = 549.5
$ws.Range("M33").Value = -170.5

$ws.Range("H34").Value = 5265.2
$ws.Range("I34").Value = 4564.154
$ws.Range("J34").Value = 6567.143
$ws.Range("K34").Value = 4564.154
$ws.Range("L34").Value = 6567.143
$ws.Range("M34").Value = -4362.154
$ws.Range("N34").Value = -6971.143

$ws.Range("H62").Value = 103488.57
$ws.Range("J62").Value = 201002.5
$ws.Range("L62").Value = 201002.5
$ws.Range("N62").Value = -202250.5

$ws.Range("H65").Value = 103488.57
$ws.Range("J65").Value = 201002.5
$ws.Range("L65").Value = 1005012.5
$ws.Range("N65").Value = -1011252.5

$ws.Range("H122").Value = 1872.5
$ws.Range("I122").Value = 1197.2
$ws.Range("K122").Value = 3591.6
$ws.Range("M122").Value = -1141.6

$ws.Range("H123").Value = 49000
$ws.Range("J123").Value = 49000
$ws.Range("L123").Value = 49000
$ws.Range("N123").Value = -58800

$ws.Range("H134").Value = 3258.9092
$ws.Range("I134").Value = 2672.4
$ws.Range("K134").Value = 8017.200000000001
$ws.Range("M134").Value = -5482.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 32068964
$ws.Range("I4").Value = 36649464
$ws.Range("J4").Value = 5471.5
$ws.Range("K4").Value = 109948392
$ws.Range("L4").Value = 16414.5
$ws.Range("M4").Value = -109948280
$ws.Range("N4").Value = -16638.5

$ws.Range("H7").Value = 10000143
$ws.Range("I7").Value = 12500158
$ws.Range("J7").Value = 84.5
$ws.Range("K7").Value = 37500474
$ws.Range("L7").Value = 253.5
$ws.Range("M7").Value = -37500362
$ws.Range("N7").Value = -477.5

$ws.Range("H10").Value = 40
$ws.Range("I10").Value = 40
$ws.Range("K10").Value = 120
$ws.Range("M10").Value = 19

$ws.Range("H11").Value = 1749.75
$ws.Range("I11").Value = 1749.75
$ws.Range("K11").Value = 5249.25
$ws.Range("M11").Value = -5109.25

$ws.Range("H37").Value = 98833.336
$ws.Range("J37").Value = 98833.336
$ws.Range("L37").Value = 296500.008
$ws.Range("N37").Value = -296724.008

$ws.Range("H81").Value = 2500

$ws.Range("H84").Value = 2500

$ws.Range("H131").Value = 1471.0656
$ws.Range("I131").Value = 907.5
$ws.Range("J131").Value = 1510.614
$ws.Range("K131").Value = 2722.5
$ws.Range("L131").Value = 4531.842000000001
$ws.Range("M131").Value = 2317.5
$ws.Range("N131").Value = -14611.842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31262

$ws.Range("H70").Value = 7750
$ws.Range("I70").Value = 7500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7230

$ws.Range("H73").Value = 7750
$ws.Range("I73").Value = 7500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6564

$ws.Range("H113").Value = 27800844
$ws.Range("J113").Value = 41330
$ws.Range("L113").Value = 41330
$ws.Range("N113").Value = -45670

$ws.Range("H132").Value = 1601
$ws.Range("I132").Value = 1947.5
$ws.Range("K132").Value = 5842.5
$ws.Range("M132").Value = -3312.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2510500.2
$ws.Range("I2").Value = 3340666.8
$ws.Range("K2").Value = 3340666.8
$ws.Range("M2").Value = -3340554.8

$ws.Range("H16").Value = 583.1
$ws.Range("I16").Value = 583.1
$ws.Range("K16").Value = 583.1
$ws.Range("M16").Value = -413.1

$ws.Range("H22").Value = 815.2
$ws.Range("I22").Value = 815.2
$ws.Range("K22").Value = 815.2
$ws.Range("M22").Value = -520.2

$ws.Range("H27").Value = 815.2
$ws.Range("I27").Value = 815.2
$ws.Range("K27").Value = 815.2
$ws.Range("M27").Value = -708.2

$ws.Range("H46").Value = 2555.0476
$ws.Range("J46").Value = 2665.4
$ws.Range("L46").Value = 2665.4
$ws.Range("N46").Value = -3041.4

$ws.Range("H47").Value = 29500
$ws.Range("I47").Value = 29000
$ws.Range("J47").Value = 30000
$ws.Range("K47").Value = 29000
$ws.Range("L47").Value = 30000
$ws.Range("M47").Value = -28510
$ws.Range("N47").Value = -30980

$ws.Range("H52").Value = 29500
$ws.Range("I52").Value = 29000
$ws.Range("J52").Value = 30000
$ws.Range("K52").Value = 29000
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -28767
$ws.Range("N52").Value = -30466

$ws.Range("H132").Value = 1999.6666
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -11058.5

$ws.Range("H136").Value = 5636.5625
$ws.Range("I136").Value = 5378.357
$ws.Range("K136").Value = 16135.071
$ws.Range("M136").Value = -13585.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2576.7273
$ws.Range("J100").Value = 3250.75
$ws.Range("L100").Value = 6501.5
$ws.Range("N100").Value = -7583.5

$ws.Range("H107").Value = 383.72726
$ws.Range("I107").Value = 383.72726
$ws.Range("K107").Value = 1151.18178
$ws.Range("M107").Value = 768.8182200000001

$ws.Range("H120").Value = 11500
$ws.Range("J120").Value = 11500
$ws.Range("L120").Value = 11500
$ws.Range("N120").Value = -21176

$ws.Range("H125").Value = 38739.2
$ws.Range("J125").Value = 38739.2
$ws.Range("L125").Value = 38739.2
$ws.Range("N125").Value = -48579.2

$ws.Range("H132").Value = 3196.7778
$ws.Range("I132").Value = 2756.1333
$ws.Range("K132").Value = 8268.3999
$ws.Range("M132").Value = -5738.3999
